# changes to recalculate subject marks
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 - SubjectAverage values recalculated (scaled by 1.5)
$ws.Range("C10").Value = 103
$ws.Range("D10").Value = 93.5
$ws.Range("E10").Value = 77
$ws.Range("F10").Value = 83.5
$ws.Range("G10").Value = 123.5
$ws.Range("H10").Value = 102.5
$ws.Range("I10").Value = 79
$ws.Range("J10").Value = 122
$ws.Range("K10").Value = 784
$ws.Range("L10").Value = 98

# Row 11 - SubjectGrades updated to reflect new averages
$ws.Range("C11").Value = "A"
$ws.Range("D11").Value = "A"
$ws.Range("E11").Value = "B+"
$ws.Range("F11").Value = "A-"
$ws.Range("G11").Value = "A"
$ws.Range("H11").Value = "A"
$ws.Range("I11").Value = "B+"
$ws.Range("J11").Value = "A"
$ws.Range("L11").Value = "A"
